$d = $word.ActiveDocument

# Find the paragraph ending "...They get copied there any way" (last bullet
# before the run of empty ListParagraph placeholders) and append two new
# bullet items after it, matching the existing numId=2 list formatting.
$rng = $d.Content
$rng.Find.Execute("They get copied there any way", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$anchorPara = $rng.Paragraphs(1)
$anchorRange = $anchorPara.Range

$anchorRange.InsertParagraphAfter()
$newPara1 = $anchorPara.Next()
$newPara1.Range.Text = "If I don’t add a Background node nothing shows after one frame, odd"

$newPara1.Range.InsertParagraphAfter()
$newPara2 = $newPara1.Next()
$newPara2.Range.Text = "If I only add a Background node the perspective stay crazy and I see a white circle"
